# Timesheet_Group5.xlsx update (sravani, 11/03/2013)
# Mark the last days of the pay period ("OFF") for the team members in
# rows 32-39 of the "FebruaryMarch 2013" sheet, then leave the selection
# where that edit was made (mirrors what a user would do interactively:
# select the range, type the value, and the active cell/selection rectangle
# is what gets stored in the sheet view).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FebruaryMarch 2013")
$ws.Activate()

# Rows 32-35: mark AP:AS as "OFF"
$ws.Range("AP32:AS35").Value = "OFF"

# Rows 36-39: mark AP:AR as "OFF" (AS stays blank)
$ws.Range("AP36:AR39").Value = "OFF"

# Reflect the final selection left behind by the edit
$ws.Range("AS32:AS35").Select() | Out-Null
